# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.750.16"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "2.465.97"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("D15").Value = "2.913.10"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "62.697.47"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "2.470.33"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "325.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +14.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "638.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.29%  "
$ws.Range("D27").Value = "2.588.21"
$ws.Range("D28").Value = "0.0₃0975"
$ws.Range("E28").Value = "  -3.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -14.57%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("E33").Value = "  -4.17%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "151.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("E40").Value = "  -4.34%  "
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("B43").Value = "BabyDogeCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D43").Value = "0.0₆0316"
$ws.Range("E43").Value = "  -12.26%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "152.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.91%  "
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.606"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0507"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  -1.56%  "
